$d = $word.ActiveDocument

# Helper: find a sub-string within a paragraph's range and split the run
# around it by toggling Bold on/off (forces the engine to materialise a
# separate <w:r> for the matched text without altering its visible
# formatting - real Word would instead mark this split with a
# <w:proofErr> pair, which is not reachable from the object model).
function Split-Run($paragraphIndex, $searchText) {
    $p = $d.Paragraphs($paragraphIndex)
    $rng = $p.Range
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

# 1. "Solution set up with Xen Engine & HGE Engine." -> split out "Xen"
Split-Run 11 "Xen"

# 2. "Add advanced map features to the game." -> whole paragraph gets strike-through
$p34 = $d.Paragraphs(34)
$p34.Range.Font.StrikeThrough = 1

# 3. "Power pellets & eating ghosts." -> whole paragraph gets themed grey colour
$p37 = $d.Paragraphs(37)
$p37.Range.Font.TextColor.ObjectThemeColor = 12

# 4. "Implement a very basic character select screen and host/join screen."
#    -> "host/join screen" becomes struck-through
$p40 = $d.Paragraphs(40)
$rng40 = $p40.Range
$rng40.Find.Execute("host/join screen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng40.Find.Found) {
    $rng40.Font.StrikeThrough = 1
}

# 5. "Integrate basic networking features into the gameplay." -> split out "gameplay"
Split-Run 45 "gameplay"

# 6. "Implement simple Ghost and PacMan AI." -> split out "PacMan"
Split-Run 49 "PacMan"

# 7. "Polish Ghost and PacMan AI." -> split out "PacMan"
Split-Run 73 "PacMan"
